$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast-error values for rows 2-11 (Q1..Q9 plus N column) after
# fixing the naive component forecaster bug.

$data = @{
    2  = @(0.2226776683031646, 0.5618053799771385, 0.6246702391763954, 0.7903608284678558, 0.7667232613880856, 46)
    3  = @(0.293438081967545,  0.6695240432108801, 0.9182272854563196, 0.9582417677477431, 0.9225146215623552, 45)
    4  = @(0.182530046921793,  0.5838921936173209, 0.6581603071260018, 0.8112707976539041, 0.7996088567086325, 44)
    5  = @(0.327844975786004,  0.7216035043336653, 1.079816604652694,  1.03914224466754,   0.9977399924889347, 43)
    6  = @(0.2304610288088512, 0.5964001741841488, 0.7431306070296015, 0.8620502346322989, 0.840742550595554,  42)
    7  = @(0.360861710175336,  0.7403235333647489, 1.140448443186462,  1.067917807317802,  1.0175868272775,    41)
    8  = @(0.2216292893798094, 0.5674340590213579, 0.6895196896333891, 0.8303732230951267, 0.8104447324343937, 40)
    9  = @(0.4304460383242469, 0.7770904355848363, 1.255534736424988,  1.120506464249532,  1.048053168048187,  39)
    10 = @(0.2885592243077175, 0.5810449934012875, 0.7354929410298064, 0.8576088508345786, 0.818446246706525,  38)
    11 = @(0.4217362328470671, 0.7352401290682791, 1.136174966075426,  1.065915083895254,  0.9924380765910319, 37)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
